# Apply the "STM32 rebrand" edits described by the commit:
#  - Slide 1 title: "... with PIC16F877A Microcontroller" ->
#    "... with STM32 Microcontroller"
#  - Slide 5 title: "Connect the keypad to the PIC16F877A" ->
#    "Keypad Pinouts"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1 - Title 3 : "Interfacing Keypad with PIC16F877A Microcontroller"
# ---------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$titleShape = $null
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $candidate = $slide1.Shapes.Item($i)
    if ($candidate.Name -eq "Title 3") {
        $titleShape = $candidate
    }
}

$titleRange = $titleShape.TextFrame.TextRange
$titleText = $titleRange.Text
$oldTail = "PIC16F877A Microcontroller"
$newTail = "STM32 Microcontroller"
$tailStart = $titleText.IndexOf($oldTail)
if ($tailStart -ge 0) {
    $tailRange = $titleRange.Characters($tailStart + 1, $oldTail.Length)
    $tailRange.Text = $newTail
}

# ---------------------------------------------------------------
# Slide 5 - Title 4 : "Connect the keypad to the PIC16F877A"
# ---------------------------------------------------------------
$slide5 = $p.Slides.Item(5)

$pinoutShape = $null
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $candidate = $slide5.Shapes.Item($i)
    if ($candidate.Name -eq "Title 4") {
        $pinoutShape = $candidate
    }
}

$pinoutRange = $pinoutShape.TextFrame.TextRange
$pinoutText = $pinoutRange.Text
$oldHeading = "Connect the keypad to the PIC16F877A"
$headingStart = $pinoutText.IndexOf($oldHeading)
if ($headingStart -ge 0) {
    $headingRange = $pinoutRange.Characters($headingStart + 1, $oldHeading.Length)
    # Replace with the first part of the new heading, keeping its run/format,
    # then append the rest as its own run (mirrors how PowerPoint splits a
    # run when only part of its text is retyped).
    $headingRange.Text = "Keypad "
    $freshRange = $pinoutShape.TextFrame.TextRange
    [void]$freshRange.InsertAfter("Pinouts")
}
